$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1
$ws.Range("A1").Value = "UDA1361TS/N1"
$ws.Range("B1").Value = "UDA1361TS/N1"
$ws.Range("C1").Value = "АЦП; Каналы: 1; 24бит; 110квыб./с; 2,4÷3,6В; SSOP16"
$ws.Range("D1").Value = 0.000071
$ws.Range("E1").Value = "https://ce8dc832c.cloudimg.io/v7/_cdn_/09/DC/00/00/0/52624_1.jpg?width=640&height=480&wat=1&wat_url=_tme-wrk_%2Ftme_new.png&wat_scale=100p&ci_sign=b7f51d47faec248dda86032605af027a8f164846"
$ws.Range("F1").Value = "//www.tme.eu/ru/details/uda1361ts_n1/preobrazovateli-a-d-mikroskhemy/nxp/uda1361ts-n1-112/"
$ws.Range("G1").Value = "{'Тип микросхемы': 'АЦП', 'Кол-во каналов': '2', 'Разрешение преобразователя': '24бит', 'Частота обновления': '110квыб./с', 'Корпус': 'SSOP16', 'Монтаж': 'SMD', 'Интерфейс': 'I2S', 'Характеристики интегральных схем': 'stereo', 'Частота': '96кГц', 'Напряжение питания': '2,4...3,6В DC'}"
$ws.Range("H1").Value = "https://www.tme.eu/Document/caa190e31c8e714c0e9d5eed22619b2e/UDA1361TS-N1-DTE.pdf"
$ws.Range("J1").Value = "Тип микросхемы АЦП"
$ws.Range("K1").Value = "NXP"
$ws.Range("L1").Value = "UDA1361TS/N1,112"

# Row 2
$ws.Range("A2").Value = "R2K150-AC01-15"

# Row 3
$ws.Range("A3").Value = "BQ27441DRZT-G1A"
$ws.Range("B3").Value = "BQ27441DRZT-G1A"
$ws.Range("C3").Value = "Система контроля; контроллер заряда аккумуляторов; SON12"
$ws.Range("D3").Value = 0.00017
$ws.Range("E3").Value = "https://ce8dc832c.cloudimg.io/v7/_cdn_/27/67/B0/00/0/751218_1.jpg?width=640&height=480&wat=1&wat_url=_tme-wrk_%2Ftme_new.png&wat_scale=100p&ci_sign=f0ec52a3f4b002933ece394ec659d28b3fdbd260"
$ws.Range("F3").Value = "//www.tme.eu/ru/details/bq27441drzt-g1a/kontrollery-batarei-i-akkumuliat-skhemy/texas-instruments/"
$ws.Range("G3").Value = "{'Тип микросхемы': 'система контроля', 'Вид микросхемы': 'контроллер заряда аккумуляторов', 'Корпус': 'SON12', 'Выходное напряжение': '4,2В', 'Интерфейс': 'I2C', 'Рабочая температура': '-40...85°C', 'Монтаж': 'SMD', 'Количество аккумуляторов': '1 x Li-Ion / Li-Po'}"
$ws.Range("I3").Value = " Интегральные монолитные схемы — регулятор заряда аккумуляторов. Корпус SON12, Выходное напряжение 4,2В, Интерфейс I2C, Рабочая температура -40...85°C, Монтаж на поверхность печатной платы, Количество аккумуляторов 1 x Li-Ion / Li-Po, предназначены для монтажа на печатную плату радиоэлектронного оборудования общепромышленного назначения."
$ws.Range("J3").Value = "Тип микросхемы система контроля"
$ws.Range("K3").Value = "TEXAS INSTRUMENTS"
$ws.Range("L3").Value = "BQ27441DRZT-G1A"

# Row 4
$ws.Range("A4").Value = "MAX7317AEE+"
$ws.Range("B4").ClearContents()
$ws.Range("C4").ClearContents()
$ws.Range("D4").ClearContents()
$ws.Range("F4").ClearContents()
$ws.Range("G4").ClearContents()
$ws.Range("I4").ClearContents()
$ws.Range("J4").ClearContents()
$ws.Range("K4").ClearContents()
$ws.Range("L4").ClearContents()

# Row 5
$ws.Range("A5").Value = "IRFR5305PBF"
$ws.Range("B5").Value = "IRFR5305PBF"
$ws.Range("C5").Value = "Транзистор: P-MOSFET; полевой; -55В; -28А; 89Вт; DPAK"
$ws.Range("D5").Value = 0.00032
$ws.Range("E5").Value = "https://ce8dc832c.cloudimg.io/v7/_cdn_/16/7E/00/00/0/59233_1.jpg?width=640&height=480&wat=1&wat_url=_tme-wrk_%2Ftme_new.png&wat_scale=100p&ci_sign=85a92da44426a86680dc5a939d6fe19e96293048"
$ws.Range("F5").Value = "//www.tme.eu/ru/details/irfr5305pbf/tranzistory-s-kanalom-p-smd/infineon-irf/"
$ws.Range("G5").Value = "{'Тип транзистора': 'P-MOSFET', 'Технология': 'HEXFET®', 'Полярность': 'полевой', 'Напряжение сток-исток': '-55В', 'Ток стока': '-28А', 'Рассеиваемая мощность': '89Вт', 'Корпус': 'DPAK', 'Напряжение затвор-исток': '±20В', 'Сопротивление в открытом состоянии': '65мОм', 'Монтаж': 'SMD', 'Заряд затвора': '42нC', 'Вид канала': 'обогащенный'}"
$ws.Range("H5").Value = "https://www.tme.eu/Document/c3ef49bf7438fe933f8717d1acbf6b87/irfr5305.pdf"
$ws.Range("I5").Value = " Кремниевые МОП-транзисторы с P-канальной структурой. Тип транзистора P-MOSFET, Технология HEXFET, Полярность полевой, Напряжение сток-исток -55В, Ток стока -28А, Рассеиваемая мощность 89Вт, Корпус DPAK, Напряжение затвор-исток +\- 20В, Сопротивление в открытом состоянии 65мОм, Монтаж на поверхность печатной платы, Заряд затвора 42нC, рабочие температуры от -40 до 85°С, предназначены для использования в радиоэлектронном оборудовании промышленного назначения."
$ws.Range("J5").Value = "Тип транзистора P-MOSFET"
$ws.Range("K5").Value = "Infineon (IRF)"
$ws.Range("L5").Value = "IRFR5305PBF"

# Row 6
$ws.Range("A6").Value = "0456020.ER"
$ws.Range("B6").Value = "0456020.ER"
$ws.Range("C6").Value = "Предохранитель: плавкая вставка; быстродействующий; 20А; 125ВAC"
$ws.Range("D6").Value = 0.0005
$ws.Range("E6").Value = "https://ce8dc832c.cloudimg.io/v7/_cdn_/0D/DD/00/00/0/56784_1.jpg?width=640&height=480&wat=1&wat_url=_tme-wrk_%2Ftme_new.png&wat_scale=100p&ci_sign=5eda84b4c4f9995ae1e3ee2216ad2332ca0cafc7"
$ws.Range("F6").Value = "//www.tme.eu/ru/details/0456020.er/predokhraniteli-smd-ostalnye/littelfuse/"
$ws.Range("G6").Value = "{'Тип предохранителя': 'плавкая вставка', 'Характеристика предохранителя': 'быстрый', 'Ток отключения': '100А', 'Номинальный ток': '20А', 'Номинальное напряжение': '125В AC', 'Монтаж': 'SMD', 'Вид предохранителя': 'керамический', 'Размер предохранителя': '10,1x3,12x3,12мм', 'Характеристика отключения': '2I<sub>n</sub>: макс 60с', 'Материал контакта': 'латунь', 'Покрытие контакта': 'посеребренные'}"
$ws.Range("H6").Value = "https://www.tme.eu/Document/a32a0db7e672d04a049bb50c11ca8186/0456020.ER.pdf"
$ws.Range("I6").Value = " Предохранители плавкие. Характеристика предохранителя быстрый, Номинальный ток 20А, Номинальное напряжение 125В переменного тока, Монтаж на поверхность печатной платы, Вид предохранителя керамический, Размер предохранителя 10,1x3,12x3,12мм, Материал контакта латунь. Предназначены для использования в радиоэлектронном оборудовании общепромышленного назначения."
$ws.Range("J6").Value = "Тип предохранителя плавкая вставка"
$ws.Range("K6").Value = "LITTELFUSE"
$ws.Range("L6").Value = "0456020.ER"

# Extend used range to row 7 (matches target dimension A1:L7)
$ws.Range("L7").Font.Bold = $false

